# Update MiniRhex hardware BOM to M3 hardware revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the (now stale) hyperlinks for the three rows whose McMaster
# part numbers changed. Find-then-delete one at a time (re-querying the
# live collection each time) so the indices of the other hyperlinks don't
# shift unexpectedly underneath us. ---
$staleAddresses = @("`$D`$6", "`$D`$11", "`$D`$13")
foreach ($addr in $staleAddresses) {
    $match = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $match = $hl
            break
        }
    }
    if ($match -ne $null) {
        $match.Delete()
    }
}

# --- Update prices and link text for the three McMaster-Carr part rows ---
# Row 6: Screws
$ws.Range("B6").Value = 6.08
$ws.Range("D6").Value = "https://www.mcmaster.com/#94500A223"

# Row 11: Nuts
$ws.Range("B11").Value = 3.36
$ws.Range("D11").Value = "https://www.mcmaster.com/#90576A102"

# Row 13: Spacer
$ws.Range("B13").Value = 4.48
$ws.Range("D13").Value = "https://www.mcmaster.com/#93657A203"

# --- Update the selected cell shown when the workbook is next opened ---
$ws.Range("B16").Select()
